$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.432.10'
$ws.Range("E2").Value = '  +4.99%  '
$ws.Range("D3").Value = '2.243.06'
$ws.Range("E3").Value = '  +3.74%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '''228.78'
$ws.Range("E5").Value = '  +0.46%  '
$ws.Range("D6").Value = '''0.637'
$ws.Range("E6").Value = '  +2.28%  '
$ws.Range("D7").Value = '''65.23'
$ws.Range("E7").Value = '  +1.32%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '''0.408'
$ws.Range("E9").Value = '  +2.61%  '
$ws.Range("D10").Value = '''0.0894'
$ws.Range("E10").Value = '  +3.93%  '
$ws.Range("E11").Value = '  +0.77%  '
$ws.Range("D12").Value = '2.569.78'
$ws.Range("E12").Value = '  +3.56%  '
$ws.Range("D13").Value = '''16.17'
$ws.Range("E13").Value = '  +0.79%  '
$ws.Range("D14").Value = '''22.40'
$ws.Range("E14").Value = '  +0.70%  '
$ws.Range("D15").Value = '''0.827'
$ws.Range("E15").Value = '  +1.56%  '
$ws.Range("E16").Value = '  +1.21%  '
$ws.Range("D17").Value = '2.245.34'
$ws.Range("E17").Value = '  +4.01%  '
$ws.Range("D18").Value = '41.243.90'
$ws.Range("E18").Value = '  +4.69%  '
$ws.Range("D19").Value = '''74.04'
$ws.Range("E19").Value = '  +3.09%  '
$ws.Range("D20").Value = '0.0₃0913'
$ws.Range("E20").Value = '  +7.00%  '
$ws.Range("D21").Value = '''6.16'
$ws.Range("E21").Value = '  +0.49%  '
$ws.Range("D22").Value = '''253.86'
$ws.Range("E22").Value = '  +9.61%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D25").Value = '''2.33'
$ws.Range("E25").Value = '  -7.18%  '
$ws.Range("D26").Value = '''9.79'
$ws.Range("E26").Value = '  +2.83%  '
$ws.Range("D27").Value = '''172.91'
$ws.Range("E27").Value = '  +0.25%  '
$ws.Range("D28").Value = '''0.146'
$ws.Range("E28").Value = '  +4.37%  '
$ws.Range("D29").Value = '''20.45'
$ws.Range("E29").Value = '  +2.79%  '
$ws.Range("E30").Value = '  +2.13%  '
$ws.Range("E31").Value = '  +5.94%  '
$ws.Range("E32").Value = '  +2.17%  '
$ws.Range("D33").Value = '''4.69'
$ws.Range("E33").Value = '  +1.68%  '
$ws.Range("E34").Value = '  +1.87%  '
$ws.Range("D35").Value = '''7.22'
$ws.Range("E35").Value = '  +1.71%  '
$ws.Range("E36").Value = '  +2.09%  '
$ws.Range("D37").Value = '''3.84'
$ws.Range("E37").Value = '  +7.83%  '
$ws.Range("D38").Value = '''2.45'
$ws.Range("E38").Value = '  +2.01%  '
$ws.Range("D39").Value = '''0.997'
$ws.Range("E39").Value = '  -0.23%  '
$ws.Range("D40").Value = '''0.000237'
$ws.Range("E40").Value = '  +53.16%  '
$ws.Range("E41").Value = '  +14.91%  '
$ws.Range("E42").Value = '  +2.16%  '
$ws.Range("D43").Value = '''8.71'
$ws.Range("E43").Value = '  +11.47%  '
$ws.Range("D44").Value = '''17.92'
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("D45").Value = '''101.65'
$ws.Range("E45").Value = '  -2.28%  '
$ws.Range("E46").Value = '  +3.53%  '
$ws.Range("D47").Value = '1.511.75'
$ws.Range("E47").Value = '  -1.77%  '
$ws.Range("D48").Value = '''0.0942'
$ws.Range("E48").Value = '  +1.88%  '
$ws.Range("E49").Value = '  +1.01%  '
$ws.Range("E50").Value = '  +1.32%  '
$ws.Range("D51").Value = '''51.73'
$ws.Range("E51").Value = '  +11.42%  '
